$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 495  # H18
$ws.Cells.Item(18, 9).Value = 490  # I18
$ws.Cells.Item(18, 10).Value = 500  # J18
$ws.Cells.Item(18, 11).Value = 490  # K18
$ws.Cells.Item(18, 12).Value = 500  # L18
$ws.Cells.Item(18, 13).Value = -206  # M18
$ws.Cells.Item(18, 14).Value = -1068  # N18
$ws.Cells.Item(43, 8).Value = 891.0833  # H43
$ws.Cells.Item(43, 9).Value = 838  # I43
$ws.Cells.Item(43, 10).Value = 929  # J43
$ws.Cells.Item(43, 11).Value = 838  # K43
$ws.Cells.Item(43, 12).Value = 929  # L43
$ws.Cells.Item(43, 13).Value = -769  # M43
$ws.Cells.Item(43, 14).Value = -1067  # N43
$ws.Cells.Item(94, 8).Value = 2600.7144  # H94
$ws.Cells.Item(94, 9).Value = 2600.7144  # I94
$ws.Cells.Item(94, 11).Value = 2600.7144  # K94
$ws.Cells.Item(94, 13).Value = -2149.7144  # M94
$ws.Cells.Item(116, 8).Value = 4300.2607  # H116
$ws.Cells.Item(116, 9).Value = 2150  # I116
$ws.Cells.Item(116, 10).Value = 5447.067  # J116
$ws.Cells.Item(116, 11).Value = 2150  # K116
$ws.Cells.Item(116, 12).Value = 5447.067  # L116
$ws.Cells.Item(116, 13).Value = 1292  # M116
$ws.Cells.Item(116, 14).Value = -12331.067  # N116
$ws.Cells.Item(129, 8).Value = 1203.4166  # H129
$ws.Cells.Item(129, 10).Value = 1371.4333  # J129
$ws.Cells.Item(129, 12).Value = 4114.2999  # L129
$ws.Cells.Item(129, 14).Value = -14114.2999  # N129
$ws.Cells.Item(138, 8).Value = 37041004  # H138
$ws.Cells.Item(138, 9).Value = 100003070  # I138
$ws.Cells.Item(138, 10).Value = 4497.4707  # J138
$ws.Cells.Item(138, 11).Value = 300009210  # K138
$ws.Cells.Item(138, 12).Value = 13492.4121  # L138
$ws.Cells.Item(138, 13).Value = -300004070  # M138
$ws.Cells.Item(138, 14).Value = -23772.4121  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4443.017  # H32
$ws.Cells.Item(32, 9).Value = 3455.423  # I32
$ws.Cells.Item(32, 11).Value = 3455.423  # K32
$ws.Cells.Item(32, 13).Value = -3168.423  # M32
$ws.Cells.Item(44, 8).Value = 0  # H44
$ws.Cells.Item(44, 10).Value = 0  # J44
$ws.Cells.Item(44, 12).Value = 0  # L44
$ws.Cells.Item(44, 14).ClearContents()  # N44
$ws.Cells.Item(80, 8).Value = 55267.5  # H80
$ws.Cells.Item(80, 10).Value = 55267.5  # J80
$ws.Cells.Item(80, 12).Value = 55267.5  # L80
$ws.Cells.Item(80, 14).Value = -57263.5  # N80
$ws.Cells.Item(82, 8).Value = 16550  # H82
$ws.Cells.Item(82, 10).Value = 18100  # J82
$ws.Cells.Item(82, 12).Value = 18100  # L82
$ws.Cells.Item(82, 14).Value = -18822  # N82
$ws.Cells.Item(83, 8).Value = 55267.5  # H83
$ws.Cells.Item(83, 10).Value = 55267.5  # J83
$ws.Cells.Item(83, 12).Value = 165802.5  # L83
$ws.Cells.Item(83, 14).Value = -175786.5  # N83
$ws.Cells.Item(85, 8).Value = 16550  # H85
$ws.Cells.Item(85, 10).Value = 18100  # J85
$ws.Cells.Item(85, 12).Value = 18100  # L85
$ws.Cells.Item(85, 14).Value = -20596  # N85
$ws.Cells.Item(97, 8).Value = 200002240  # H97
$ws.Cells.Item(97, 9).Value = 3557  # I97
$ws.Cells.Item(97, 10).Value = 500000260  # J97
$ws.Cells.Item(97, 11).Value = 3557  # K97
$ws.Cells.Item(97, 12).Value = 500000260  # L97
$ws.Cells.Item(97, 13).Value = -3061  # M97
$ws.Cells.Item(97, 14).Value = -500001252  # N97
$ws.Cells.Item(140, 8).Value = 40429  # H140
$ws.Cells.Item(140, 10).Value = 40429  # J140
$ws.Cells.Item(140, 12).Value = 40429  # L140
$ws.Cells.Item(140, 14).Value = -50789  # N140

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 726.0454999999999  # H80
$ws.Cells.Item(80, 10).Value = 780.3333  # J80
$ws.Cells.Item(80, 12).Value = 780.3333  # L80
$ws.Cells.Item(80, 14).Value = -2776.3333  # N80
$ws.Cells.Item(83, 8).Value = 726.0454999999999  # H83
$ws.Cells.Item(83, 10).Value = 780.3333  # J83
$ws.Cells.Item(83, 12).Value = 3901.6665  # L83
$ws.Cells.Item(83, 14).Value = -13885.6665  # N83
$ws.Cells.Item(99, 8).Value = 1104.3158  # H99
$ws.Cells.Item(99, 9).Value = 1113.3334  # I99
$ws.Cells.Item(99, 11).Value = 1113.3334  # K99
$ws.Cells.Item(99, 13).Value = 384.6666  # M99
$ws.Cells.Item(105, 8).Value = 4549387.5  # H105
$ws.Cells.Item(105, 9).Value = 5026.6665  # I105
$ws.Cells.Item(105, 10).Value = 10002620  # J105
$ws.Cells.Item(105, 11).Value = 5026.6665  # K105
$ws.Cells.Item(105, 12).Value = 10002620  # L105
$ws.Cells.Item(105, 13).Value = -3279.6665  # M105
$ws.Cells.Item(105, 14).Value = -10006114  # N105
$ws.Cells.Item(107, 8).Value = 802.5714  # H107
$ws.Cells.Item(107, 9).Value = 802.5714  # I107
$ws.Cells.Item(107, 10).Value = 0  # J107
$ws.Cells.Item(107, 11).Value = 802.5714  # K107
$ws.Cells.Item(107, 12).Value = 0  # L107
$ws.Cells.Item(107, 13).Value = 1117.4286  # M107
$ws.Cells.Item(107, 14).ClearContents()  # N107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 789  # H16
$ws.Cells.Item(16, 9).Value = 735  # I16
$ws.Cells.Item(16, 10).Value = 897  # J16
$ws.Cells.Item(16, 11).Value = 735  # K16
$ws.Cells.Item(16, 12).Value = 897  # L16
$ws.Cells.Item(16, 13).Value = -448  # M16
$ws.Cells.Item(16, 14).Value = -1471  # N16
$ws.Cells.Item(31, 8).Value = 2871.5  # H31
$ws.Cells.Item(31, 10).Value = 3278.8276  # J31
$ws.Cells.Item(31, 12).Value = 3278.8276  # L31
$ws.Cells.Item(31, 14).Value = -3868.8276  # N31
$ws.Cells.Item(34, 8).Value = 2871.5  # H34
$ws.Cells.Item(34, 10).Value = 3278.8276  # J34
$ws.Cells.Item(34, 12).Value = 3278.8276  # L34
$ws.Cells.Item(34, 14).Value = -3682.8276  # N34
$ws.Cells.Item(99, 8).Value = 23813300  # H99
$ws.Cells.Item(99, 10).Value = 71433416  # J99
$ws.Cells.Item(99, 12).Value = 71433416  # L99
$ws.Cells.Item(99, 14).Value = -71436412  # N99
$ws.Cells.Item(113, 8).Value = 789  # H113
$ws.Cells.Item(113, 9).Value = 735  # I113
$ws.Cells.Item(113, 10).Value = 897  # J113
$ws.Cells.Item(113, 11).Value = 735  # K113
$ws.Cells.Item(113, 12).Value = 897  # L113
$ws.Cells.Item(113, 13).Value = 1435  # M113
$ws.Cells.Item(113, 14).Value = -5237  # N113
$ws.Cells.Item(122, 8).Value = 1592.2222  # H122
$ws.Cells.Item(122, 10).Value = 1993.3334  # J122
$ws.Cells.Item(122, 12).Value = 5980.0002  # L122
$ws.Cells.Item(122, 14).Value = -10880.0002  # N122
$ws.Cells.Item(126, 8).Value = 23813300  # H126
$ws.Cells.Item(126, 10).Value = 71433416  # J126
$ws.Cells.Item(126, 12).Value = 214300248  # L126
$ws.Cells.Item(126, 14).Value = -214305188  # N126

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 38462390  # H86
$ws.Cells.Item(86, 10).Value = 83334330  # J86
$ws.Cells.Item(86, 12).Value = 250002990  # L86
$ws.Cells.Item(86, 14).Value = -250005362  # N86
$ws.Cells.Item(89, 8).Value = 38462390  # H89
$ws.Cells.Item(89, 10).Value = 83334330  # J89
$ws.Cells.Item(89, 12).Value = 750008970  # L89
$ws.Cells.Item(89, 14).Value = -750020826  # N89
$ws.Cells.Item(114, 8).Value = 1274.8572  # H114
$ws.Cells.Item(114, 10).Value = 1843.5555  # J114
$ws.Cells.Item(114, 12).Value = 5530.666499999999  # L114
$ws.Cells.Item(114, 14).Value = -12038.6665  # N114
$ws.Cells.Item(122, 8).Value = 529.1667  # H122
$ws.Cells.Item(122, 9).Value = 340.6  # I122
$ws.Cells.Item(122, 10).Value = 663.8570999999999  # J122
$ws.Cells.Item(122, 11).Value = 3065.4  # K122
$ws.Cells.Item(122, 12).Value = 5974.7139  # L122
$ws.Cells.Item(122, 13).Value = -615.4000000000001  # M122
$ws.Cells.Item(122, 14).Value = -10874.7139  # N122
$ws.Cells.Item(131, 8).Value = 789.6  # H131
$ws.Cells.Item(131, 10).Value = 823.59784  # J131
$ws.Cells.Item(131, 12).Value = 2470.79352  # L131
$ws.Cells.Item(131, 14).Value = -12550.79352  # N131
$ws.Cells.Item(140, 8).Value = 1500.5454  # H140
$ws.Cells.Item(140, 9).Value = 639.38464  # I140
$ws.Cells.Item(140, 10).Value = 2744.4443  # J140
$ws.Cells.Item(140, 11).Value = 1918.15392  # K140
$ws.Cells.Item(140, 12).Value = 8233.332900000001  # L140
$ws.Cells.Item(140, 13).Value = 3261.84608  # M140
$ws.Cells.Item(140, 14).Value = -18593.3329  # N140
$ws.Cells.Item(141, 8).Value = 1853.75  # H141
$ws.Cells.Item(141, 9).Value = 1853.75  # I141
$ws.Cells.Item(141, 11).Value = 5561.25  # K141
$ws.Cells.Item(141, 13).Value = -381.25  # M141

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(96, 8).Value = 0  # H96
$ws.Cells.Item(96, 10).Value = 0  # J96
$ws.Cells.Item(96, 12).Value = 0  # L96
$ws.Cells.Item(96, 14).ClearContents()  # N96
$ws.Cells.Item(97, 8).Value = 1403.05  # H97
$ws.Cells.Item(97, 9).Value = 1473.3334  # I97
$ws.Cells.Item(97, 11).Value = 1473.3334  # K97
$ws.Cells.Item(97, 13).Value = -977.3334  # M97
$ws.Cells.Item(113, 8).Value = 2153.5  # H113
$ws.Cells.Item(113, 9).Value = 1710.5333  # I113
$ws.Cells.Item(113, 10).Value = 2891.7778  # J113
$ws.Cells.Item(113, 11).Value = 1710.5333  # K113
$ws.Cells.Item(113, 12).Value = 2891.7778  # L113
$ws.Cells.Item(113, 13).Value = 459.4666999999999  # M113
$ws.Cells.Item(113, 14).Value = -7231.7778  # N113
$ws.Cells.Item(132, 8).Value = 31869.588  # H132
$ws.Cells.Item(132, 9).Value = 1987.1  # I132
$ws.Cells.Item(132, 10).Value = 74558.86  # J132
$ws.Cells.Item(132, 11).Value = 5961.299999999999  # K132
$ws.Cells.Item(132, 12).Value = 223676.58  # L132
$ws.Cells.Item(132, 13).Value = -3431.299999999999  # M132
$ws.Cells.Item(132, 14).Value = -228736.58  # N132
$ws.Cells.Item(136, 8).Value = 17591.834  # H136
$ws.Cells.Item(136, 10).Value = 17591.834  # J136
$ws.Cells.Item(136, 12).Value = 52775.50199999999  # L136
$ws.Cells.Item(136, 14).Value = -57875.50199999999  # N136

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 4212.375  # H22
$ws.Cells.Item(22, 10).Value = 3849.8333  # J22
$ws.Cells.Item(22, 12).Value = 3849.8333  # L22
$ws.Cells.Item(22, 14).Value = -4439.8333  # N22
$ws.Cells.Item(27, 8).Value = 4212.375  # H27
$ws.Cells.Item(27, 10).Value = 3849.8333  # J27
$ws.Cells.Item(27, 12).Value = 3849.8333  # L27
$ws.Cells.Item(27, 14).Value = -4063.8333  # N27
$ws.Cells.Item(61, 8).Value = 4274.375  # H61
$ws.Cells.Item(61, 9).Value = 1563.1818  # I61
$ws.Cells.Item(61, 10).Value = 10239  # J61
$ws.Cells.Item(61, 11).Value = 1563.1818  # K61
$ws.Cells.Item(61, 12).Value = 10239  # L61
$ws.Cells.Item(61, 13).Value = -1361.1818  # M61
$ws.Cells.Item(61, 14).Value = -10643  # N61
$ws.Cells.Item(113, 8).Value = 4274.375  # H113
$ws.Cells.Item(113, 9).Value = 1563.1818  # I113
$ws.Cells.Item(113, 10).Value = 10239  # J113
$ws.Cells.Item(113, 11).Value = 1563.1818  # K113
$ws.Cells.Item(113, 12).Value = 10239  # L113
$ws.Cells.Item(113, 13).Value = 606.8181999999999  # M113
$ws.Cells.Item(113, 14).Value = -14579  # N113
$ws.Cells.Item(132, 8).Value = 929456.25  # H132
$ws.Cells.Item(132, 9).Value = 1508129.1  # I132
$ws.Cells.Item(132, 10).Value = 3579.6  # J132
$ws.Cells.Item(132, 11).Value = 4524387.300000001  # K132
$ws.Cells.Item(132, 12).Value = 10738.8  # L132
$ws.Cells.Item(132, 13).Value = -4521857.300000001  # M132
$ws.Cells.Item(132, 14).Value = -15798.8  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 13000  # H32
$ws.Cells.Item(32, 10).Value = 13000  # J32
$ws.Cells.Item(32, 12).Value = 13000  # L32
$ws.Cells.Item(32, 14).Value = -13634  # N32
$ws.Cells.Item(81, 8).Value = 66667760  # H81
$ws.Cells.Item(81, 9).Value = 1239.1  # I81
$ws.Cells.Item(81, 11).Value = 2478.2  # K81
$ws.Cells.Item(81, 13).Value = -1417.2  # M81
$ws.Cells.Item(84, 8).Value = 66667760  # H84
$ws.Cells.Item(84, 9).Value = 1239.1  # I84
$ws.Cells.Item(84, 11).Value = 12391  # K84
$ws.Cells.Item(84, 13).Value = -7087  # M84
$ws.Cells.Item(113, 8).Value = 2457203  # H113
$ws.Cells.Item(113, 9).Value = 214  # I113
$ws.Cells.Item(113, 11).Value = 642  # K113
$ws.Cells.Item(113, 13).Value = 1528  # M113
$ws.Cells.Item(136, 8).Value = 20241672  # H136
$ws.Cells.Item(136, 9).Value = 27165770  # I136
$ws.Cells.Item(136, 10).Value = 2000.3846  # J136
$ws.Cells.Item(136, 11).Value = 81497310  # K136
$ws.Cells.Item(136, 12).Value = 6001.1538  # L136
$ws.Cells.Item(136, 13).Value = -81494760  # M136
$ws.Cells.Item(136, 14).Value = -11101.1538  # N136
